$wb = $excel.ActiveWorkbook

# --- "Basis" sheet: move several user stories to their finished status ---
$basis = $wb.Worksheets.Item("Basis")

$basis.Range("A13").Value = "Afgerond"

$basis.Range("A19").Value = "Afgerond"
$basis.Range("A20").Value = "Afgerond"
$basis.Range("A21").Value = "Afgerond"

$basis.Range("A24").Value = "Afgerond"
$basis.Range("A25").Value = "Kom er niet uit"
$basis.Range("A26").Value = "Afgerond"
$basis.Range("A27").Value = "Afgerond"
$basis.Range("A28").Value = "Afgerond"
$basis.Range("A29").Value = "Afgerond"
$basis.Range("A30").Value = "Afgerond"

# Reflect where the author was last working: scrolled down to row 27,
# with C28 selected.
$basis.Activate()
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$basis.Range("C28").Select()

# --- "Extra" sheet: same status update on its single data row ---
$extra = $wb.Worksheets.Item("Extra")
$extra.Range("A2").Value = "Afgerond"
$extra.Range("A4").Select()

# Leave "Basis" as the selected/active sheet, matching tabSelected="1"
$basis.Activate()
